# feat: add 2022-Q3 data
#
# Before: sheet "总计" (totals) + sheet "2022-Q2" (fund-holding detail for Q2).
# After : sheet "总计" gets a new row for Q3 (and the old Q2 row shifts down);
#         the fund-holding sheet is repurposed to hold the NEW Q3 numbers and
#         renamed "2022-Q3"; a fresh copy of the original sheet (still holding
#         the old Q2 numbers) is inserted right after it and named "2022-Q2",
#         preserving the historical snapshot.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# 1) Duplicate the Q2 detail sheet so the old data survives, placed
#    immediately after the source sheet. Rename the original to its new
#    quarter first so the copy can take over the "2022-Q2" name.
$q2Sheet.Copy($null, $q2Sheet)
$q2Sheet.Name = "2022-Q3"
$q2Copy = $wb.Worksheets.Item($q2Sheet.Index + 1)
$q2Copy.Name = "2022-Q2"

# 2) Turn the original sheet into the new "2022-Q3" detail sheet by
#    overwriting its values in place (already renamed above).
$q2Sheet.Range("B2").Value = "'003318"
$q2Sheet.Range("B2").Style = "Normal"
$q2Sheet.Range("C2").Value = "景顺长城中证500行业中性低波动指数"
$q2Sheet.Range("D2").Value = "'10.25"
$q2Sheet.Range("D2").Style = "Normal"
$q2Sheet.Range("E2").Value = "'93.67"
$q2Sheet.Range("E2").Style = "Normal"
$q2Sheet.Range("F2").Value = "'1.14"
$q2Sheet.Range("F2").Style = "Normal"
$q2Sheet.Range("G2").Value = "'0.1168"
$q2Sheet.Range("G2").Style = "Normal"
$q2Sheet.Range("H2").Value = 5

$q2Sheet.Range("B3").Value = "'512260"
$q2Sheet.Range("B3").Style = "Normal"
$q2Sheet.Range("C3").Value = "华安中证500行业中性低波动ETF"
$q2Sheet.Range("D3").Value = "'1.07"
$q2Sheet.Range("D3").Style = "Normal"
$q2Sheet.Range("E3").Value = "'97.91"
$q2Sheet.Range("E3").Style = "Normal"
$q2Sheet.Range("F3").Value = "'1.20"
$q2Sheet.Range("F3").Style = "Normal"
$q2Sheet.Range("G3").Value = "'0.0128"
$q2Sheet.Range("G3").Style = "Normal"
$q2Sheet.Range("H3").Value = 5

# 3) Update the totals sheet: the existing Q2 row's numbers move to a new
#    row 3, and row 2 becomes the new Q3 totals. A3 carries the same
#    "index column" style as A2, so copy that formatting down first.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.15

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("D2").Value = 0.13
